$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B/C shift to C/D)
$ws.Columns.Item(2).Insert()

# New header for inserted column B
$ws.Range("B1").Value = "StatQuery"

# New "stat" Cypher query text for inserted column B, row 2
$statQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Miniature Schnauzer']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$ws.Range("B2").Value = $statQuery

# Match formatting of column A (wrap text) on the new B2 cell
$ws.Range("B2").WrapText = $true

# Match width of column A for the new column B
$ws.Columns.Item(2).ColumnWidth = 75

# Reset selection to A2 (single cell) as in the final saved view
$ws.Range("A2").Select()
